# Updated multiple management files to have consistent and english values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Replace fertilizer abbreviation "KAS" with "CAN" in the Notice column (E)
# for the 1st/2nd/3rd nitrogen application rows across all four treatment blocks.
$canCells = @("E47","E48","E58","E59","E69","E70","E71","E80","E81","E82")
foreach ($cellRef in $canCells) {
    $ws.Range($cellRef).Value = "CAN"
}

# 2) Clear the now-removed "others" / "Kieserit" / "Bittersalz" fertilization
# entries (amount, date, notice) for each of the four treatment blocks.
$rowsToClear = @(50, 51, 61, 62, 72, 73, 83, 84)
foreach ($r in $rowsToClear) {
    $ws.Range("B" + $r).Value = ""
    $ws.Range("C" + $r).Value = ""
    $ws.Range("D" + $r).Value = ""
    $ws.Range("E" + $r).Value = ""
}

$wb.Save()
